{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" block\n// (an empty paragraph, the \"Ver no Jupiter...\" paragraph, another empty\n// paragraph, and the page-break paragraph that immediately preceded it),\n// leaving the rest of the document (including the final empty paragraph and\n// the trailing page-break paragraph) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst marker = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  // Paragraph layout around the marker (as found in the source document):\n  //   [markerIndex - 1] empty paragraph\n  //   [markerIndex]     \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  //   [markerIndex + 1] empty paragraph\n  //   [markerIndex + 2] empty paragraph with a page break before it\n  const toDelete = [];\n  if (markerIndex - 1 >= 0 && items[markerIndex - 1].text === \"\") {\n    toDelete.push(items[markerIndex - 1]);\n  }\n  toDelete.push(items[markerIndex]);\n  if (markerIndex + 1 < items.length && items[markerIndex + 1].text === \"\") {\n    toDelete.push(items[markerIndex + 1]);\n  }\n  if (markerIndex + 2 < items.length && items[markerIndex + 2].text === \"\") {\n    toDelete.push(items[markerIndex + 2]);\n  }\n\n  // Delete from last to first so earlier indices stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph holding the \"Ver no Jupiter Salvar em pdf Salvar em\n# docx\" text. That paragraph is flanked by an empty paragraph right before\n# it and, right after it, another empty paragraph followed by an empty\n# page-break paragraph \u2014 all four need to go, leaving the rest of the\n# document (including the final empty paragraph and its trailing page-break\n# paragraph) untouched.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $beforeEmpty = $target.Previous()\n    $afterEmpty1 = $target.Next()\n    $afterEmpty2 = $afterEmpty1.Next()\n\n    $startPos = $beforeEmpty.Range.Start\n    $endPos = $afterEmpty2.Range.End\n\n    $deleteRange = $d.Range($startPos, $endPos)\n    $deleteRange.Delete()\n}\n"}
